# Web Developement Angela Yu Section 12 update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the shared string used as a note, now referenced by I9 instead of I4
$ws.Range("I4").ClearContents()
$ws.Range("I9").Value = "Refer each time you create a new Website"

# Fill in Completion Date (column C) for rows 5-13
$ws.Range("C5").Value = 44055
$ws.Range("C6").Value = 44055
$ws.Range("C7").Value = 44056
$ws.Range("C8").Value = 44057
$ws.Range("C9").Value = 44057
$ws.Range("C10").Value = 44062
$ws.Range("C11").Value = 44062
$ws.Range("C12").Value = 44062
$ws.Range("C13").Value = 44062
$ws.Range("C5:C13").NumberFormat = "m/d/yy"

# Update selection / view to match
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("C25").Select()
